$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (header values)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON)
$ws.Range("B2").Value = 826.86838530406374
$ws.Range("C2").Value = 439.3877776586998
$ws.Range("D2").Value = 939.63583694062663
$ws.Range("E2").Value = 417.24887536970186

# Row 3 (STR)
$ws.Range("B3").Value = 770.17302114195684
$ws.Range("C3").Value = 444.87390981478114
$ws.Range("D3").Value = 853.00131378641288
$ws.Range("E3").Value = 359.55529361567272

# Update selection to match the new data extent used in the diff
$ws.Range("B1:E3").Select()
